# CORE_holdings.xlsx update
# - Bump the "as of" date in the confidential disclosure note from 2021-03-22 to 2021-03-23
# - Refresh the Weight (col D) and Percent Change (col E) figures for rows 2-8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet ships protected; temporarily unprotect so the cells can be updated,
# then restore protection afterwards.
$ws.Unprotect()

# Update the confidential disclosure text (date changed from 2021-03-22 to 2021-03-23)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# Row 2 (IVE)
$ws.Range("D2").Value = 0.4963000646746448
$ws.Range("E2").Value = -0.01273386924663877

# Row 3 (IVW)
$ws.Range("D3").Value = 0.2436669400226923
$ws.Range("E3").Value = -0.002785946447918186

# Row 4 (IJK)
$ws.Range("D4").Value = 0.09865250457777726
$ws.Range("E4").Value = -0.02280297207276438

# Row 5 (IJJ)
$ws.Range("D5").Value = 0.1020872417628093
$ws.Range("E5").Value = -0.02967184122617406

# Row 6 (IJS)
$ws.Range("D6").Value = 0.03065325403333874
$ws.Range("E6").Value = -0.03720427996466091

# Row 7 (IJT)
$ws.Range("D7").Value = 0.02863999492873752
$ws.Range("E7").Value = -0.03220587102242067

# Row 8 (Total)
$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = -0.01434015813993805

# Restore sheet protection (original password is unknown to this automation,
# so re-protect without one to keep the sheet locked for normal editing).
$ws.Protect()
